# Apply updated bus voltage magnitude (vm_pu) results for the 380 kV case (B2=1.02)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.026388868250889
$ws.Range("D2").Value = 1.027797155090833
$ws.Range("E2").Value = 1.026606404435739
$ws.Range("F2").Value = 1.036174583677398
$ws.Range("I2").Value = 1.031582048339116
$ws.Range("J2").Value = 1.031552584253658
$ws.Range("K2").Value = 1.030615837652987
$ws.Range("L2").Value = 1.029428559723848
$ws.Range("M2").Value = 1.038969079276475
$ws.Range("N2").Value = 1.014486322699843

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.02751292824
$ws.Range("D3").Value = 1.028776941328492
$ws.Range("E3").Value = 1.027565808511736
$ws.Range("F3").Value = 1.037471758884877
$ws.Range("I3").Value = 1.031774013593289
$ws.Range("J3").Value = 1.032315484495395
$ws.Range("K3").Value = 1.031402960026777
$ws.Range("L3").Value = 1.030195102367082
$ws.Range("M3").Value = 1.040074507691915
$ws.Range("N3").Value = 1.014744112155595

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.028240057962907
$ws.Range("D4").Value = 1.029411021598052
$ws.Range("E4").Value = 1.028186774746934
$ws.Range("F4").Value = 1.038311015686693
$ws.Range("I4").Value = 1.031896632981637
$ws.Range("J4").Value = 1.032808405131871
$ws.Range("K4").Value = 1.031911764147895
$ws.Range("L4").Value = 1.030690660207686
$ws.Range("M4").Value = 1.040789150077184
$ws.Range("N4").Value = 1.014910539157136

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.02854569373961
$ws.Range("D5").Value = 1.029677611668997
$ws.Range("E5").Value = 1.028447869441967
$ws.Range("F5").Value = 1.038663816780895
$ws.Range("I5").Value = 1.031947800407854
$ws.Range("J5").Value = 1.033015455498957
$ws.Range("K5").Value = 1.032125542181326
$ws.Range("L5").Value = 1.0308988862448
$ws.Range("M5").Value = 1.041089433280557
$ws.Range("N5").Value = 1.014980414181722

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.028597008461709
$ws.Range("D6").Value = 1.029722374660561
$ws.Range("E6").Value = 1.028491710774749
$ws.Range("F6").Value = 1.03872305237847
$ws.Range("I6").Value = 1.031956369259121
$ws.Range("J6").Value = 1.033050210009942
$ws.Range("K6").Value = 1.032161429242001
$ws.Range("L6").Value = 1.030933842080716
$ws.Range("M6").Value = 1.041139843259635
$ws.Range("N6").Value = 1.014992141178734

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.028244142076278
$ws.Range("D7").Value = 1.02941458369778
$ws.Range("E7").Value = 1.028190263347494
$ws.Range("F7").Value = 1.038315729917932
$ws.Range("I7").Value = 1.031897318183275
$ws.Range("J7").Value = 1.03281117242996
$ws.Range("K7").Value = 1.03191462114402
$ws.Range("L7").Value = 1.030693442952064
$ws.Range("M7").Value = 1.040793163073037
$ws.Range("N7").Value = 1.014911473187649

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.026768794291647
$ws.Range("D8").Value = 1.02812825943977
$ws.Range("E8").Value = 1.026930605038508
$ws.Range("F8").Value = 1.03661299251281
$ws.Range("I8").Value = 1.031647254189144
$ws.Range("J8").Value = 1.031810560782295
$ws.Range("K8").Value = 1.030881956616558
$ws.Range("L8").Value = 1.029687709102597
$ws.Range("M8").Value = 1.039342798266077
$ws.Range("N8").Value = 1.014573522592237

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.024167360766162
$ws.Range("D9").Value = 1.025862283766746
$ws.Range("E9").Value = 1.024712193538241
$ws.Range("F9").Value = 1.03361167717236
$ws.Range("I9").Value = 1.031194390309051
$ws.Range("J9").Value = 1.03004176395247
$ws.Range("K9").Value = 1.029058290641968
$ws.Range("L9").Value = 1.027912038683891
$ws.Range("M9").Value = 1.036782065700433
$ws.Range("N9").Value = 1.013975095750661

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.022431858898439
$ws.Range("D10").Value = 1.024352066678945
$ws.Range("E10").Value = 1.023234078498005
$ws.Range("F10").Value = 1.031610082553061
$ws.Range("I10").Value = 1.030884259160341
$ws.Range("J10").Value = 1.028858769936395
$ws.Range("K10").Value = 1.027839804763226
$ws.Range("L10").Value = 1.026725916973165
$ws.Range("M10").Value = 1.035071439458959
$ws.Range("N10").Value = 1.013574177834437

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.021680060718178
$ws.Range("D11").Value = 1.023698220762973
$ws.Range("E11").Value = 1.022594227525481
$ws.Range("F11").Value = 1.030743170077919
$ws.Range("I11").Value = 1.030748016961818
$ws.Range("J11").Value = 1.028345611080178
$ws.Range("K11").Value = 1.027311536295867
$ws.Range("L11").Value = 1.026211750504537
$ws.Range("M11").Value = 1.034329871605453
$ws.Range("N11").Value = 1.01340010759205

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.021400760417644
$ws.Range("D12").Value = 1.023455365709233
$ws.Range("E12").Value = 1.022356585133458
$ws.Range("F12").Value = 1.03042112615726
$ws.Range("I12").Value = 1.030697116846028
$ws.Range("J12").Value = 1.02815486274507
$ws.Range("K12").Value = 1.027115214500114
$ws.Range("L12").Value = 1.026020680133337
$ws.Range("M12").Value = 1.03405428955744
$ws.Range("N12").Value = 1.013335379317686

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.021460673524845
$ws.Range("D13").Value = 1.023507458392931
$ws.Range("E13").Value = 1.022407559033864
$ws.Range("F13").Value = 1.03049020725406
$ws.Range("I13").Value = 1.030708048387713
$ws.Range("J13").Value = 1.028195785182334
$ws.Range("K13").Value = 1.027157330701195
$ws.Range("L13").Value = 1.026061669282267
$ws.Range("M13").Value = 1.034113408793609
$ws.Range("N13").Value = 1.013349266961016

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.021656974677545
$ws.Range("D14").Value = 1.023678146039491
$ws.Range("E14").Value = 1.022574583389539
$ws.Range("F14").Value = 1.030716550536543
$ws.Range("I14").Value = 1.030743815535015
$ws.Range("J14").Value = 1.028329846589385
$ws.Range("K14").Value = 1.027295310306683
$ws.Range("L14").Value = 1.026195958327684
$ws.Range("M14").Value = 1.034307094587058
$ws.Range("N14").Value = 1.01339475857891

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.021777915753552
$ws.Range("D15").Value = 1.023783313955374
$ws.Range("E15").Value = 1.02267749613145
$ws.Range("F15").Value = 1.030856003514442
$ws.Range("I15").Value = 1.030765813930424
$ws.Range("J15").Value = 1.028412427896218
$ws.Range("K15").Value = 1.027380310904485
$ws.Range("L15").Value = 1.026278686816242
$ws.Range("M15").Value = 1.034426413407643
$ws.Range("N15").Value = 1.013422778073563

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.022481746494412
$ws.Range("D16").Value = 1.024395462066241
$ws.Range("E16").Value = 1.023276547088534
$ws.Range("F16").Value = 1.031667612053047
$ws.Range("I16").Value = 1.030893259931948
$ws.Range("J16").Value = 1.028892807254793
$ws.Range("K16").Value = 1.02787485031113
$ws.Range("L16").Value = 1.026760028480708
$ws.Range("M16").Value = 1.035120636663893
$ws.Range("N16").Value = 1.013585720380592

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.022923155514632
$ws.Range("D17").Value = 1.024779469416634
$ws.Range("E17").Value = 1.023652364331996
$ws.Range("F17").Value = 1.032176654872238
$ws.Range("I17").Value = 1.030972680311572
$ws.Range("J17").Value = 1.029193891045251
$ws.Range("K17").Value = 1.028184885609555
$ws.Range("L17").Value = 1.027061808826315
$ws.Range("M17").Value = 1.035555873995378
$ws.Range("N17").Value = 1.013687803773039

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.023180591708506
$ws.Range("D18").Value = 1.025003462956444
$ws.Range("E18").Value = 1.02387158983569
$ws.Range("F18").Value = 1.032473550969264
$ws.Range("I18").Value = 1.031018816419241
$ws.Range("J18").Value = 1.029369419935573
$ws.Range("K18").Value = 1.02836566070945
$ws.Range("L18").Value = 1.027237777447056
$ws.Range("M18").Value = 1.035809658182144
$ws.Range("N18").Value = 1.013747301940882

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.023268365725103
$ws.Range("D19").Value = 1.025079840439235
$ws.Range("E19").Value = 1.023946343083977
$ws.Range("F19").Value = 1.032574781559872
$ws.Range("I19").Value = 1.031034515681122
$ws.Range("J19").Value = 1.029429255818907
$ws.Range("K19").Value = 1.028427289643888
$ws.Range("L19").Value = 1.027297768934859
$ws.Range("M19").Value = 1.035896178124788
$ws.Range("N19").Value = 1.013767581592895

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.022875799642309
$ws.Range("D20").Value = 1.024738268183423
$ws.Range("E20").Value = 1.023612040876173
$ws.Range("F20").Value = 1.032122041479416
$ws.Range("I20").Value = 1.030964178750922
$ws.Range("J20").Value = 1.029161596732069
$ws.Range("K20").Value = 1.028151628303457
$ws.Range("L20").Value = 1.02702943631928
$ws.Range("M20").Value = 1.035509185685911
$ws.Range("N20").Value = 1.013676855883083

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.021599170246944
$ws.Range("D21").Value = 1.023627882445099
$ws.Range("E21").Value = 1.022525398146897
$ws.Range("F21").Value = 1.03064989904204
$ws.Range("I21").Value = 1.030733291110909
$ws.Range("J21").Value = 1.028290372670849
$ws.Range("K21").Value = 1.027254681501375
$ws.Range("L21").Value = 1.026156415932162
$ws.Range("M21").Value = 1.034250062574952
$ws.Range("N21").Value = 1.013381364388345

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.02079621841177
$ws.Range("D22").Value = 1.022929810389239
$ws.Range("E22").Value = 1.021842337150197
$ws.Range("F22").Value = 1.029724106464692
$ws.Range("I22").Value = 1.03058642353372
$ws.Range("J22").Value = 1.027741799000117
$ws.Range("K22").Value = 1.026690160145991
$ws.Range("L22").Value = 1.025607015134913
$ws.Range("M22").Value = 1.033457644442417
$ws.Range("N22").Value = 1.013195167221052

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.02122190590094
$ws.Range("D23").Value = 1.023299865081485
$ws.Range("E23").Value = 1.022204426301726
$ws.Range("F23").Value = 1.030214906303414
$ws.Range("I23").Value = 1.030664441970395
$ws.Range("J23").Value = 1.02803268445577
$ws.Range("K23").Value = 1.026989478399083
$ws.Range("L23").Value = 1.025898310362415
$ws.Range("M23").Value = 1.033877792827202
$ws.Range("N23").Value = 1.01329391278801

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.022897197815481
$ws.Range("D24").Value = 1.024756885220512
$ws.Range("E24").Value = 1.023630261255307
$ws.Range("F24").Value = 1.03214671898182
$ws.Range("I24").Value = 1.030968020822423
$ws.Range("J24").Value = 1.029176189414297
$ws.Range("K24").Value = 1.028166656043666
$ws.Range("L24").Value = 1.027044064230649
$ws.Range("M24").Value = 1.035530282377232
$ws.Range("N24").Value = 1.013681802903338

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.024840101450961
$ws.Range("D25").Value = 1.026448014050105
$ws.Range("E25").Value = 1.025285558566381
$ws.Range("F25").Value = 1.03438770610058
$ws.Range("I25").Value = 1.031312915164177
$ws.Range("J25").Value = 1.030499706703261
$ws.Range("K25").Value = 1.029530226958876
$ws.Range("L25").Value = 1.028371502575133
$ws.Range("M25").Value = 1.014130149553705
